$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 / column F (ChildId description): new description text, Times New Roman 11
# dark-gray font color (#1F1F1F), centered horizontally.
$ws.Range("F5").Value = "Primary id for auto increatment."
$ws.Range("F5").Font.Name = "Times New Roman"
$ws.Range("F5").Font.Size = 11
$ws.Range("F5").Font.Color = 2039583
$ws.Range("F5").HorizontalAlignment = -4108

# Row 6 / column B (AlbumId column name cell): highlight with green fill.
$ws.Range("B6").Interior.Color = 5296274

# Row 6 / column F (AlbumId description): new description text, Times New Roman 11,
# dark-gray font color (#1F1F1F).
$ws.Range("F6").Value = "Get Albumid  (Album) in Numaric "
$ws.Range("F6").Font.Name = "Times New Roman"
$ws.Range("F6").Font.Size = 11
$ws.Range("F6").Font.Color = 2039583

# Row 7 / column F (ImgPath description): new description text, default formatting.
$ws.Range("F7").Value = "Get ImagePath Comma Separated"

# Row 8 / column F (IsActive description): new description text, Times New Roman 11,
# dark-gray font color (#1F1F1F).
$ws.Range("F8").Value = "For active or inactive record."
$ws.Range("F8").Font.Name = "Times New Roman"
$ws.Range("F8").Font.Size = 11
$ws.Range("F8").Font.Color = 2039583

# Widen column F to fit the new description text.
$ws.Columns(6).ColumnWidth = 37.17

# Move the active selection.
[void]$ws.Range("C12").Select()
